# Auto-generated edit script: updates cryptos list D (Price) and E (Volume 1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.904.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.549.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.487"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.22"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.771.18"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.547.68"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.900.86"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.69"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.48"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0697"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.416.84"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.04%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.526"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.684.65"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.20%  "
